# "Generate Report for Handback"
#
# The localization status workbook is updated to reflect that handback has
# completed and is in sync with en-US: the status text flips from
# "Ready for handoff" to "Handed back: in sync with en-US" everywhere it is
# used (Overview + per-language sheets), the per-language "Latest Handback
# DateTime" timestamps are stamped with the real handback time (instead of
# the zero-date placeholder), and a "Latest Target File" / "Latest Handback
# File" pair of hyperlinked cells is populated for each row on the
# per-language sheets.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/8216ed58e901c988c5612056f495ca34e07797f9/e2e"
$row1Md = "14d7b9cd-bc5d-4809-a087-4877fc173bc5.md"
$row2Md = "f3289838-938f-400a-a1ff-af6edfd3fca0.md"

function Style-AsLink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = 15570276
}

# ---------------------------------------------------------------------
# Overview sheet: refresh the status text (shared by both language columns)
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = $statusNew
$ovw.Range("C2").Value = $statusNew
$ovw.Range("B3").Value = $statusNew
$ovw.Range("C3").Value = $statusNew

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew

$zh.Range("H2").Value = "2016-03-14 06:35:18"
$zh.Range("H3").Value = "2016-03-14 06:35:18"

# Row 2 - Latest Target File / Latest Handback File
$zh.Range("F2").Value = $row1Md
$zh.Hyperlinks.Add($zh.Range("F2"), "$mdBase/$row1Md", "", "", $row1Md)
Style-AsLink $zh.Range("F2")

$zhG2Name = "14d7b9cd-bc5d-4809-a087-4877fc173bc5.91b60eaee5d2a3c8345595d58b726b1eb48a2583.zh-cn.xlf"
$zh.Range("G2").Value = $zhG2Name
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1768a8dc9b42b2beb35e0f86cf664312a0cfc15d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/14d7b9cd-bc5d-4809-a087-4877fc173bc5.91b60eaee5d2a3c8345595d58b726b1eb48a2583.zh-cn.xlf", "", "", $zhG2Name)
Style-AsLink $zh.Range("G2")

# Row 3 - Latest Target File / Latest Handback File
$zh.Range("F3").Value = $row2Md
$zh.Hyperlinks.Add($zh.Range("F3"), "$mdBase/$row2Md", "", "", $row2Md)
Style-AsLink $zh.Range("F3")

$zhG3Name = "f3289838-938f-400a-a1ff-af6edfd3fca0.2d13511705896270bd09b80c3a360fc98dfc1e19.zh-cn.xlf"
$zh.Range("G3").Value = $zhG3Name
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1768a8dc9b42b2beb35e0f86cf664312a0cfc15d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f3289838-938f-400a-a1ff-af6edfd3fca0.2d13511705896270bd09b80c3a360fc98dfc1e19.zh-cn.xlf", "", "", $zhG3Name)
Style-AsLink $zh.Range("G3")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew

$de.Range("H2").Value = "2016-03-14 06:35:24"
$de.Range("H3").Value = "2016-03-14 06:35:24"

# Row 2 - Latest Target File / Latest Handback File
$de.Range("F2").Value = $row1Md
$de.Hyperlinks.Add($de.Range("F2"), "$mdBase/$row1Md", "", "", $row1Md)
Style-AsLink $de.Range("F2")

$deG2Name = "14d7b9cd-bc5d-4809-a087-4877fc173bc5.91b60eaee5d2a3c8345595d58b726b1eb48a2583.de-de.xlf"
$de.Range("G2").Value = $deG2Name
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02ca1539b8b075bea92dc121454bedc70b3b4549/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/14d7b9cd-bc5d-4809-a087-4877fc173bc5.91b60eaee5d2a3c8345595d58b726b1eb48a2583.de-de.xlf", "", "", $deG2Name)
Style-AsLink $de.Range("G2")

# Row 3 - Latest Target File / Latest Handback File
$de.Range("F3").Value = $row2Md
$de.Hyperlinks.Add($de.Range("F3"), "$mdBase/$row2Md", "", "", $row2Md)
Style-AsLink $de.Range("F3")

$deG3Name = "f3289838-938f-400a-a1ff-af6edfd3fca0.2d13511705896270bd09b80c3a360fc98dfc1e19.de-de.xlf"
$de.Range("G3").Value = $deG3Name
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02ca1539b8b075bea92dc121454bedc70b3b4549/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f3289838-938f-400a-a1ff-af6edfd3fca0.2d13511705896270bd09b80c3a360fc98dfc1e19.de-de.xlf", "", "", $deG3Name)
Style-AsLink $de.Range("G3")
